# Update crypto price/volume table values per the Mar 20 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. '335.38', '27.899.22' with
# two dots, or values with significant trailing zeros like '0.08710'/'1.440').
# They must stay as literal text (as in the source file), so we prefix each with
# a leading apostrophe -- exactly like a user typing '335.38 into a cell -- which
# forces Excel to store it as text instead of auto-converting it to a number.

$ws.Range('D2').Value = "'" + '27.899.22'
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').Value = "'" + '1.741.32'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'" + '335.38'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').Value = "'" + '1.004'
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').Value = "'" + '0.3757'
$ws.Range('E7').Value = '  -3.85%  '
$ws.Range('D8').Value = "'" + '0.3341'
$ws.Range('E8').Value = '  -4.18%  '
$ws.Range('D9').Value = "'" + '44.89'
$ws.Range('E9').Value = '  -7.06%  '
$ws.Range('D10').Value = "'" + '1.108'
$ws.Range('E10').Value = '  -5.89%  '
$ws.Range('D11').Value = "'" + '0.07189'
$ws.Range('E11').Value = '  -4.55%  '
$ws.Range('E12').Value = '  +0.21%  '
$ws.Range('D13').Value = "'" + '22.29'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = "'" + '6.127'
$ws.Range('E14').Value = '  -5.96%  '
$ws.Range('D15').Value = "'" + '7.093'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').Value = "'" + '1.744.43'
$ws.Range('E16').Value = '  -3.77%  '
$ws.Range('D17').Value = "'" + '0.00001053'
$ws.Range('E17').Value = '  -4.55%  '
$ws.Range('E18').Value = '  -2.27%  '
$ws.Range('D19').Value = "'" + '78.79'
$ws.Range('E19').Value = '  -7.49%  '
$ws.Range('D20').Value = "'" + '1.002'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').Value = "'" + '16.75'
$ws.Range('E21').Value = '  -5.48%  '
$ws.Range('D22').Value = "'" + '6.214'
$ws.Range('E22').Value = '  -5.09%  '
$ws.Range('D23').Value = "'" + '27.894.59'
$ws.Range('E23').Value = '  -1.13%  '
$ws.Range('D24').Value = "'" + '11.59'
$ws.Range('E24').Value = '  -6.99%  '
$ws.Range('D25').Value = "'" + '2.395'
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = "'" + '152.86'
$ws.Range('E26').Value = '  -1.28%  '
$ws.Range('D27').Value = "'" + '19.64'
$ws.Range('E27').Value = '  -8.07%  '
$ws.Range('D28').Value = "'" + '2.305'
$ws.Range('E28').Value = '  -8.29%  '
$ws.Range('D29').Value = "'" + '1.943.23'
$ws.Range('E29').Value = '  -3.75%  '
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = "'" + '1.257'
$ws.Range('E30').Value = '  -14.44%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = "'" + '131.05'
$ws.Range('E31').Value = '  -4.21%  '
$ws.Range('D32').Value = "'" + '4.033'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('D33').Value = "'" + '5.752'
$ws.Range('E33').Value = '  -9.31%  '
$ws.Range('D34').Value = "'" + '0.08710'
$ws.Range('E34').Value = '  -1.48%  '
$ws.Range('D35').Value = "'" + '12.17'
$ws.Range('E35').Value = '  -7.52%  '
$ws.Range('D36').Value = "'" + '0.6633'
$ws.Range('E36').Value = '  -4.00%  '
$ws.Range('D37').Value = "'" + '0.02311'
$ws.Range('E37').Value = '  -6.50%  '
$ws.Range('D38').Value = "'" + '0.06213'
$ws.Range('E38').Value = '  -4.98%  '
$ws.Range('D39').Value = "'" + '5.142'
$ws.Range('E39').Value = '  -6.01%  '
$ws.Range('D40').Value = "'" + '0.2107'
$ws.Range('E40').Value = '  -5.12%  '
$ws.Range('D41').Value = "'" + '1.212'
$ws.Range('E41').Value = '  -4.12%  '
$ws.Range('D42').Value = "'" + '1.440'
$ws.Range('E42').Value = '  -11.34%  '
$ws.Range('D43').Value = "'" + '1.003'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = "'" + '7.938'
$ws.Range('E44').Value = '  -6.96%  '
$ws.Range('D45').Value = "'" + '13.73'
$ws.Range('E45').Value = '  -6.00%  '
$ws.Range('D46').Value = "'" + '3.813'
$ws.Range('E46').Value = '  -1.61%  '
$ws.Range('D47').Value = "'" + '0.6004'
$ws.Range('E47').Value = '  -6.35%  '
$ws.Range('D48').Value = "'" + '127.02'
$ws.Range('E48').Value = '  -3.74%  '
$ws.Range('D49').Value = "'" + '2.009'
$ws.Range('E49').Value = '  -7.02%  '
$ws.Range('D50').Value = "'" + '48.84'
$ws.Range('E50').Value = '  +10.36%  '
$ws.Range('D51').Value = "'" + '0.07097'
$ws.Range('E51').Value = '  -1.84%  '
